$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 32551.406
$ws.Range("I132").Value = 34658
$ws.Range("J132").Value = 952.5
$ws.Range("K132").Value = 103974
$ws.Range("L132").Value = 2857.5
$ws.Range("M132").Value = -101444
$ws.Range("N132").Value = -7917.5

$ws.Range("H133").Value = 42995.383
$ws.Range("J133").Value = 42995.383
$ws.Range("L133").Value = 42995.383
$ws.Range("N133").Value = -53115.383

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18479.072
$ws.Range("I32").Value = 4366.3823
$ws.Range("K32").Value = 4366.3823
$ws.Range("M32").Value = -4079.3823

$ws.Range("H52").Value = 37000
$ws.Range("J52").Value = 37000
$ws.Range("L52").Value = 37000
$ws.Range("N52").Value = -37636

$ws.Range("H61").Value = 3786.0232
$ws.Range("I61").Value = 3205.0557
$ws.Range("K61").Value = 3205.0557
$ws.Range("M61").Value = -2993.0557

$ws.Range("H122").Value = 2556.6553
$ws.Range("I122").Value = 1306.7894
$ws.Range("K122").Value = 3920.3682
$ws.Range("M122").Value = -1470.3682

$ws.Range("H132").Value = 2520.6135
$ws.Range("I132").Value = 2006.8379
$ws.Range("J132").Value = 5236.2856
$ws.Range("K132").Value = 6020.5137
$ws.Range("L132").Value = 15708.8568
$ws.Range("M132").Value = -3490.5137
$ws.Range("N132").Value = -20768.8568

$ws.Range("H136").Value = 3786.0232
$ws.Range("I136").Value = 3205.0557
$ws.Range("K136").Value = 9615.167099999999
$ws.Range("M136").Value = -7065.167099999999

$ws.Range("H139").Value = 48020.715
$ws.Range("J139").Value = 48020.715
$ws.Range("L139").Value = 48020.715
$ws.Range("N139").Value = -58300.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2126.6667
$ws.Range("I134").Value = 1189.7021
$ws.Range("J134").Value = 5514.154
$ws.Range("K134").Value = 3569.1063
$ws.Range("L134").Value = 16542.462
$ws.Range("M134").Value = -1034.1063
$ws.Range("N134").Value = -21612.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 919.8
$ws.Range("I16").Value = 880.46155
$ws.Range("J16").Value = 992.8570999999999
$ws.Range("K16").Value = 880.46155
$ws.Range("L16").Value = 992.8570999999999
$ws.Range("M16").Value = -593.46155
$ws.Range("N16").Value = -1566.8571

$ws.Range("H20").Value = 43919.6
$ws.Range("J20").Value = 43919.6
$ws.Range("L20").Value = 43919.6
$ws.Range("N20").Value = -44391.6

$ws.Range("H30").Value = 43919.6
$ws.Range("J30").Value = 43919.6
$ws.Range("L30").Value = 43919.6
$ws.Range("N30").Value = -44101.6

$ws.Range("H99").Value = 5669.6816
$ws.Range("I99").Value = 8457.546
$ws.Range("J99").Value = 2881.818
$ws.Range("K99").Value = 8457.546
$ws.Range("L99").Value = 2881.818
$ws.Range("M99").Value = -6959.546
$ws.Range("N99").Value = -5877.818

$ws.Range("H107").Value = 457.5909
$ws.Range("I107").Value = 298.26666
$ws.Range("J107").Value = 799
$ws.Range("K107").Value = 298.26666
$ws.Range("L107").Value = 799
$ws.Range("M107").Value = 1621.73334
$ws.Range("N107").Value = -4639

$ws.Range("H113").Value = 919.8
$ws.Range("I113").Value = 880.46155
$ws.Range("J113").Value = 992.8570999999999
$ws.Range("K113").Value = 880.46155
$ws.Range("L113").Value = 992.8570999999999
$ws.Range("M113").Value = 1289.53845
$ws.Range("N113").Value = -5332.8571

$ws.Range("H126").Value = 5669.6816
$ws.Range("I126").Value = 8457.546
$ws.Range("J126").Value = 2881.818
$ws.Range("K126").Value = 25372.638
$ws.Range("L126").Value = 8645.454000000002
$ws.Range("M126").Value = -22902.638
$ws.Range("N126").Value = -13585.454

$ws.Range("H128").Value = 43919.6
$ws.Range("J128").Value = 43919.6
$ws.Range("L128").Value = 43919.6
$ws.Range("N128").Value = -53879.6

$ws.Range("H132").Value = 3706.8333
$ws.Range("I132").Value = 3320.2173
$ws.Range("J132").Value = 4977.143
$ws.Range("K132").Value = 9960.651899999999
$ws.Range("L132").Value = 14931.429
$ws.Range("M132").Value = -7430.651899999999
$ws.Range("N132").Value = -19991.429

$ws.Range("H134").Value = 2445.2666
$ws.Range("I134").Value = 1130.7916
$ws.Range("J134").Value = 7703.1665
$ws.Range("K134").Value = 3392.3748
$ws.Range("L134").Value = 23109.4995
$ws.Range("M134").Value = -857.3748000000001
$ws.Range("N134").Value = -28179.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1006.1579
$ws.Range("I97").Value = 801.2143
$ws.Range("J97").Value = 1580
$ws.Range("K97").Value = 801.2143
$ws.Range("L97").Value = 1580
$ws.Range("M97").Value = -305.2143
$ws.Range("N97").Value = -2572

$ws.Range("H122").Value = 2166
$ws.Range("I122").Value = 2110.182
$ws.Range("J122").Value = 2780
$ws.Range("K122").Value = 6330.545999999999
$ws.Range("L122").Value = 8340
$ws.Range("M122").Value = -3880.545999999999
$ws.Range("N122").Value = -13240

$ws.Range("H132").Value = 3111.75
$ws.Range("I132").Value = 2676.923
$ws.Range("J132").Value = 4996
$ws.Range("K132").Value = 8030.768999999999
$ws.Range("L132").Value = 14988
$ws.Range("M132").Value = -5500.768999999999
$ws.Range("N132").Value = -20048

$ws.Range("H137").Value = 43000
$ws.Range("J137").Value = 43000
$ws.Range("L137").Value = 43000
$ws.Range("N137").Value = -53200

$ws.Range("H139").Value = 41279.715
$ws.Range("J139").Value = 41279.715
$ws.Range("L139").Value = 41279.715
$ws.Range("N139").Value = -51559.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2563.9285
$ws.Range("I40").Value = 1524.25
$ws.Range("K40").Value = 1524.25
$ws.Range("M40").Value = -1388.25

$ws.Range("H61").Value = 2893.2307
$ws.Range("I61").Value = 2952
$ws.Range("J61").Value = 2882.5454
$ws.Range("K61").Value = 2952
$ws.Range("L61").Value = 2882.5454
$ws.Range("M61").Value = -2750
$ws.Range("N61").Value = -3286.5454

$ws.Range("H93").Value = 1922.9524
$ws.Range("I93").Value = 1765.3077
$ws.Range("J93").Value = 2179.125
$ws.Range("K93").Value = 1765.3077
$ws.Range("L93").Value = 2179.125
$ws.Range("M93").Value = -517.3077000000001
$ws.Range("N93").Value = -4675.125

$ws.Range("H113").Value = 2893.2307
$ws.Range("I113").Value = 2952
$ws.Range("J113").Value = 2882.5454
$ws.Range("K113").Value = 2952
$ws.Range("L113").Value = 2882.5454
$ws.Range("M113").Value = -782
$ws.Range("N113").Value = -7222.5454

$ws.Range("H136").Value = 3009.389
$ws.Range("I136").Value = 1694.4333
$ws.Range("K136").Value = 5083.2999
$ws.Range("M136").Value = -2533.2999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2620.6155
$ws.Range("I81").Value = 1335.5333
$ws.Range("J81").Value = 4373
$ws.Range("K81").Value = 2671.0666
$ws.Range("L81").Value = 8746
$ws.Range("M81").Value = -1610.0666
$ws.Range("N81").Value = -10868

$ws.Range("H84").Value = 2620.6155
$ws.Range("I84").Value = 1335.5333
$ws.Range("J84").Value = 4373
$ws.Range("K84").Value = 13355.333
$ws.Range("L84").Value = 43730
$ws.Range("M84").Value = -8051.333000000001
$ws.Range("N84").Value = -54338

$ws.Range("H122").Value = 1117.3334
$ws.Range("I122").Value = 991.63635
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 2974.90905
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -524.9090500000002
$ws.Range("N122").Value = -12400

$ws.Range("H126").Value = 40817
$ws.Range("I126").Value = 77348.46000000001
$ws.Range("J126").Value = 1241.25
$ws.Range("K126").Value = 232045.38
$ws.Range("L126").Value = 3723.75
$ws.Range("M126").Value = -229575.38
$ws.Range("N126").Value = -8663.75

$ws.Range("H132").Value = 1750.28
$ws.Range("I132").Value = 1747.6875
$ws.Range("J132").Value = 1760.65
$ws.Range("K132").Value = 5243.0625
$ws.Range("L132").Value = 5281.950000000001
$ws.Range("M132").Value = -2713.0625
$ws.Range("N132").Value = -10341.95

$ws.Range("H136").Value = 1051.2545
$ws.Range("I136").Value = 549.475
$ws.Range("J136").Value = 2389.3333
$ws.Range("K136").Value = 1648.425
$ws.Range("L136").Value = 7167.999899999999
$ws.Range("M136").Value = 901.5749999999998
$ws.Range("N136").Value = -12267.9999
